$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for Wins/Losses/Ties, matching the style of the existing header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 71
    $ws.Cells.Item($r, 32).Value = 0
}
